$p = $ppt.ActivePresentation

# Locate the slide by its persistent SlideID (262) rather than a raw
# positional index, since that is how the change is addressed.
$target = $null
foreach ($sl in $p.Slides) {
    if ($sl.SlideID -eq 262) {
        $target = $sl
        break
    }
}

if ($target -ne $null) {
    # The picture (shape id 5124, "Picture 4") was nudged slightly to the
    # right/down, e.g. after the "vibration" figure was re-aligned for the
    # start of the waves material. Size stays the same, only position moves.
    foreach ($shp in $target.Shapes) {
        if ($shp.Id -eq 5124) {
            # PowerPoint's Shape.Left/Top are expressed in points while the
            # OOXML stores EMUs (1 pt = 12700 EMU). The literal point values
            # below are chosen so that PowerPoint's internal single-precision
            # round trip reproduces the exact target EMU offsets:
            #   Left -> 892016 EMU, Top -> 2159797 EMU
            $shp.Left = 70.2375221
            $shp.Top = 170.0627559
        }
    }
}
